$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Task 9 "Sound" (row 22): mark as done, clear its note
$ws.Range("J22").Value = $true
$ws.Range("G22").ClearContents()

# Task 13 "Level Select Menu & Multiple Levels" (row 26): mark as done, clear its note
$ws.Range("J26").Value = $true
$ws.Range("G26").ClearContents()

# Task 14 "Speed Run Mode" (row 27): add a note (stays "To Be Done")
$ws.Range("G27").Value = "y"

# Task 17 "World Tilt Mode" (row 30): mark as done, clear its note
$ws.Range("J30").Value = $true
$ws.Range("G30").ClearContents()
